$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -5,6)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: -9,-9)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -10,1)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: 7,2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 2,-2)"
$ws.Range("F1").Value = "(305251175, Or  Leder: 3,-9)"

$ws.Range("A3").Value = "cost: 596.4148259012275"
$ws.Range("A4").Value = "time: 71.42685323765343"
